$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("s137","s137_e101_24019-32_1_5.jpeg","meltpatch","1140","846","25","25","72","1"),
    @("s138","s138_e111_24019-32_2_1.jpeg","meltpatch","1098","269","51","25","89","2"),
    @("s139","s139_e104_24019-32_1_8.jpeg","meltpatch","340","920","77","25","77","3"),
    @("s140","s140_e113_24019-32_3_1.jpeg","meltpatch","1372","1556","103","25","103","4"),
    @("s141","s141_e107_24019-32_2_5.jpeg","meltpatch","2020","916","129","25","5","5"),
    @("s142","s142_e114_24019-32_3_2.jpeg","meltpatch","2521","851","155","25","9","6"),
    @("s143","s143_e106_24019-32_2_6.jpeg","meltpatch","2302","1604","181","25","171","7"),
    @("s144","s144_e112_24019-32_2_0.jpeg","meltpatch","339","1314","207","25","21","8"),
    @("s145","s145_e99_24019-32_1_3.jpeg","meltpatch","1506","788","233","25","100","9"),
    @("s146","s146_e103_24019-32_1_7.jpeg","meltpatch","876","764","259","25","56","10"),
    @("s147","s147_e100_24019-32_1_4.jpeg","meltpatch","515","668","38","38","29","1"),
    @("s148","s148_e122_24019-32_4_6.jpeg","meltpatch","892","1274","77","38","116","2"),
    @("s149","s149_e117_24019-32_3_5.jpeg","meltpatch","1823","1698","116","38","57","3"),
    @("s150","s150_e119_24019-32_3_7.jpeg","meltpatch","922","400","155","38","14","4"),
    @("s151","s151_e105_24019-32_2_7.jpeg","meltpatch","118","780","194","38","127","5"),
    @("s152","s152_e121_24019-32_4_7.jpeg","meltpatch","2606","1792","233","38","150","6"),
    @("s153","s153_e129_24019-32_1_5.jpeg","meltpatch","51","928","19","19","13","1"),
    @("s154","s154_e130_24019-32_1_6.jpeg","meltpatch","460","525","39","19","26","2"),
    @("s155","s155_e145_24019-32_3_5.jpeg","meltpatch","296","757","58","19","41","3"),
    @("s156","s156_e127_24019-32_1_3.jpeg","meltpatch","1439","455","78","19","146","4"),
    @("s157","s157_e144_24019-32_3_4.jpeg","meltpatch","796","492","98","19","128","5"),
    @("s158","s158_e125_24019-32_1_1.jpeg","meltpatch","1982","763","117","19","29","6"),
    @("s159","s159_e148_24019-32_3_8.jpeg","meltpatch","651","766","137","19","3","7"),
    @("s160","s160_e143_24019-32_3_3.jpeg","meltpatch","909","32","157","19","119","8"),
    @("s161","s161_e135_24019-32_2_5.jpeg","meltpatch","1910","575","176","19","148","9"),
    @("s162","s162_e134_24019-32_2_6.jpeg","meltpatch","235","129","196","19","15","10"),
    @("s163","s163_e132_24019-32_1_8.jpeg","meltpatch","857","238","29","29","24","1"),
    @("s164","s164_e136_24019-32_2_4.jpeg","meltpatch","1085","1247","58","29","107","2"),
    @("s165","s165_e139_24019-32_2_1.jpeg","meltpatch","944","560","88","29","125","3"),
    @("s166","s166_e152_24019-32_4_4.jpeg","meltpatch","1890","854","117","29","121","4"),
    @("s167","s167_e126_24019-32_1_2.jpeg","meltpatch","574","726","147","29","8","5"),
    @("s168","s168_e142_24019-32_3_2.jpeg","meltpatch","1414","1124","176","29","37","6")
)

$startRow = 138
$numericCols = @(4,5,6,7,8,9)  # D,E,F,G,H,I (1-indexed columns within the 9-col row)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $rows[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $val = $rowVals[$c - 1]
        if ($numericCols -contains $c) {
            # Force text storage for numeric-looking values (matches source data
            # where every cell, including numbers, is stored as text).
            $ws.Cells.Item($r, $c).Value = "'" + $val
        } else {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}

# Update the sheet's used-range dimension and view state to reflect the newly
# appended rows (138-169), mirroring how Excel itself updates these on edit.
$excel.ActiveWindow.ScrollRow = 124
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 55
[void]$ws.Range("A138:I169").Select()
